# Adds a new "2022-Q3" quarter sheet (right after the "总计" summary sheet)
# and records its totals (6 holdings, 1.03亿元) as the new first data row of
# the "总计" sheet, pushing the older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the 2022-Q3 row at the top of the
#    data and shift the existing quarters down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 1.03

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 8
$summary.Range("D3").Value = 1.74

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 2.28

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 15
$summary.Range("D5").Value = 6.34

# Row 6 is brand new - seed its index cell's formatting from the row above
# (A2:A5 already carry the bold/border/center style) before filling values.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 7
$summary.Range("D6").Value = 0.85

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q3" worksheet, inserted right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Columns B:G hold values that *look* numeric (fund codes with leading
# zeros, percentages with trailing zeros, ...) but must round-trip as
# plain text, exactly like every other quarter sheet in this workbook.
# Pre-format as Text so the leading/trailing zeros survive the assignment.
$q3.Range("B2:G7").NumberFormat = "@"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @(0, "001766", "上投摩根医疗健康股票A",               "8.99", "82.99", "5.25", "0.4720", 6),
    @(1, "002666", "前海开源沪港深创新成长灵活配置混合A", "6.14", "89.54", "4.69", "0.2880", 4),
    @(2, "630010", "华商价值精选混合",                     "4.30", "81.81", "3.09", "0.1329", 9),
    @(3, "002667", "前海开源沪港深创新成长灵活配置混合C", "2.27", "89.54", "4.69", "0.1065", 4),
    @(4, "630006", "华商产业升级混合",                     "0.85", "81.97", "3.12", "0.0265", 9),
    @(5, "014932", "上投摩根医疗健康股票C",               "0.04", "82.99", "5.25", "0.0021", 6)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Reuse the exact header / index-column styles already defined in the
# workbook (bold + thin border + centered) instead of re-deriving new
# cellXfs via ad-hoc Font/Border property writes.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

# Strip the temporary "@" text format back off B:G (values are already
# committed as text, so this just restores the plain/default cell style
# used everywhere else in the sheet).
$summary.Range("C2").Copy()
$q3.Range("B2:G7").PasteSpecial(-4122)

$excel.CutCopyMode = 0
